# Allocated Tasks to team members in sheet
# - Fill in the Sprint number on the title cell.
# - Replace the "Team Member" placeholder in the "Assigned To" column
#   with the actual team members the tasks were allocated to.
# - Update the active selection on the Burn Down Chart sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burn Down Chart")

# Fill in the Sprint # in the title/header cell (A1, merged A1:G1)
$title = $ws.Range("A1").Value2
$ws.Range("A1").Value = $title + "1"

# Assign tasks to team members (column E, "Assigned To") for each User Story row
$ws.Range("E4").Value  = "Mark Pratt"
$ws.Range("E5").Value  = "Mark Pratt"
$ws.Range("E6").Value  = "Patrick Garcia"
$ws.Range("E7").Value  = "Patrick Garcia"
$ws.Range("E8").Value  = "Patrick Garcia"
$ws.Range("E9").Value  = "Mark Pratt"
$ws.Range("E10").Value = "Patrick Garcia"
$ws.Range("E11").Value = "Mark Pratt"

# Update the saved selection/active cell on the sheet
$ws.Activate()
$ws.Range("A11").Select()
